$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 153.70589
$ws.Range("I55").Value = 35.77778
$ws.Range("J55").Value = 286.375
$ws.Range("K55").Value = 35.77778
$ws.Range("L55").Value = 286.375
$ws.Range("M55").Value = 178.22222
$ws.Range("N55").Value = -714.375
$ws.Range("H58").Value = 3269624
$ws.Range("I58").Value = 6536248
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 19608744
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -19608594
$ws.Range("N58").Value = -9300
$ws.Range("H116").Value = 2333.3333
$ws.Range("I116").Value = 2333.3333
$ws.Range("K116").Value = 2333.3333
$ws.Range("M116").Value = 1108.6667
$ws.Range("H129").Value = 3945
$ws.Range("J129").Value = 1019.12
$ws.Range("L129").Value = 3057.36
$ws.Range("N129").Value = -13057.36
$ws.Range("H138").Value = 3181.3914
$ws.Range("I138").Value = 1935.762
$ws.Range("J138").Value = 4227.72
$ws.Range("K138").Value = 5807.286
$ws.Range("L138").Value = 12683.16
$ws.Range("M138").Value = -667.2860000000001
$ws.Range("N138").Value = -22963.16
$ws.Range("H141").Value = 2519.5908
$ws.Range("I141").Value = 2407.875
$ws.Range("J141").Value = 2817.5
$ws.Range("K141").Value = 7223.625
$ws.Range("L141").Value = 8452.5
$ws.Range("M141").Value = -2043.625
$ws.Range("N141").Value = -18812.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32029.402
$ws.Range("I32").Value = 9596.064
$ws.Range("K32").Value = 9596.064
$ws.Range("M32").Value = -9309.064
$ws.Range("H61").Value = 1148.6
$ws.Range("I61").Value = 1204.3158
$ws.Range("K61").Value = 1204.3158
$ws.Range("M61").Value = -992.3158000000001
$ws.Range("H74").Value = 989.3333
$ws.Range("I74").Value = 924.6667
$ws.Range("J74").Value = 1118.6666
$ws.Range("K74").Value = 924.6667
$ws.Range("L74").Value = 1118.6666
$ws.Range("M74").Value = -50.66669999999999
$ws.Range("N74").Value = -2866.6666
$ws.Range("H77").Value = 989.3333
$ws.Range("I77").Value = 924.6667
$ws.Range("J77").Value = 1118.6666
$ws.Range("K77").Value = 4623.3335
$ws.Range("L77").Value = 5593.333000000001
$ws.Range("M77").Value = -255.3334999999997
$ws.Range("N77").Value = -14329.333
$ws.Range("H132").Value = 15294.488
$ws.Range("I132").Value = 16728.053
$ws.Range("K132").Value = 50184.159
$ws.Range("M132").Value = -47654.159
$ws.Range("H136").Value = 1148.6
$ws.Range("I136").Value = 1204.3158
$ws.Range("K136").Value = 3612.9474
$ws.Range("M136").Value = -1062.9474

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1953.171
$ws.Range("I134").Value = 1819.3235
$ws.Range("J134").Value = 3090.875
$ws.Range("K134").Value = 5457.970499999999
$ws.Range("L134").Value = 9272.625
$ws.Range("M134").Value = -2922.970499999999
$ws.Range("N134").Value = -14342.625

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40139.184
$ws.Range("I31").Value = 856.1
$ws.Range("J31").Value = 54168.855
$ws.Range("K31").Value = 856.1
$ws.Range("L31").Value = 54168.855
$ws.Range("M31").Value = -561.1
$ws.Range("N31").Value = -54758.855
$ws.Range("H34").Value = 40139.184
$ws.Range("I34").Value = 856.1
$ws.Range("J34").Value = 54168.855
$ws.Range("K34").Value = 856.1
$ws.Range("L34").Value = 54168.855
$ws.Range("M34").Value = -654.1
$ws.Range("N34").Value = -54572.855
$ws.Range("H132").Value = 2973.2173
$ws.Range("I132").Value = 2884.3
$ws.Range("J132").Value = 3566
$ws.Range("K132").Value = 8652.900000000001
$ws.Range("L132").Value = 10698
$ws.Range("M132").Value = -6122.900000000001
$ws.Range("N132").Value = -15758

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 503.94116
$ws.Range("I23").Value = 25.5
$ws.Range("J23").Value = 567.73334
$ws.Range("K23").Value = 76.5
$ws.Range("L23").Value = 1703.20002
$ws.Range("M23").Value = 158.5
$ws.Range("N23").Value = -2173.20002
$ws.Range("H131").Value = 830582.7
$ws.Range("J131").Value = 1013172.7
$ws.Range("L131").Value = 3039518.1
$ws.Range("N131").Value = -3049598.1

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2783.3572
$ws.Range("I122").Value = 2108.7778
$ws.Range("J122").Value = 3997.6
$ws.Range("K122").Value = 6326.3334
$ws.Range("L122").Value = 11992.8
$ws.Range("M122").Value = -3876.3334
$ws.Range("N122").Value = -16892.8
$ws.Range("H126").Value = 2403.3044
$ws.Range("I126").Value = 2359.7896
$ws.Range("J126").Value = 2610
$ws.Range("K126").Value = 7079.3688
$ws.Range("L126").Value = 7830
$ws.Range("M126").Value = -4609.3688
$ws.Range("N126").Value = -12770
$ws.Range("H132").Value = 4716.3335
$ws.Range("I132").Value = 3663.5454
$ws.Range("J132").Value = 6370.7144
$ws.Range("K132").Value = 10990.6362
$ws.Range("L132").Value = 19112.1432
$ws.Range("M132").Value = -8460.636200000001
$ws.Range("N132").Value = -24172.1432
$ws.Range("H135").Value = 44710.223
$ws.Range("J135").Value = 44710.223
$ws.Range("L135").Value = 44710.223
$ws.Range("N135").Value = -54850.223

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1476.8
$ws.Range("I22").Value = 3556.6667
$ws.Range("J22").Value = 585.4286
$ws.Range("K22").Value = 3556.6667
$ws.Range("L22").Value = 585.4286
$ws.Range("M22").Value = -3261.6667
$ws.Range("N22").Value = -1175.4286
$ws.Range("H27").Value = 1476.8
$ws.Range("I27").Value = 3556.6667
$ws.Range("J27").Value = 585.4286
$ws.Range("K27").Value = 3556.6667
$ws.Range("L27").Value = 585.4286
$ws.Range("M27").Value = -3449.6667
$ws.Range("N27").Value = -799.4286
$ws.Range("H93").Value = 2500.28
$ws.Range("I93").Value = 2769.5334
$ws.Range("J93").Value = 2096.4
$ws.Range("K93").Value = 2769.5334
$ws.Range("L93").Value = 2096.4
$ws.Range("M93").Value = -1521.5334
$ws.Range("N93").Value = -4592.4
$ws.Range("H132").Value = 5919.4614
$ws.Range("I132").Value = 6772.222
$ws.Range("J132").Value = 4000.75
$ws.Range("K132").Value = 20316.666
$ws.Range("L132").Value = 12002.25
$ws.Range("M132").Value = -17786.666
$ws.Range("N132").Value = -17062.25

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 15840
$ws.Range("J32").Value = 15840
$ws.Range("L32").Value = 15840
$ws.Range("N32").Value = -16474
$ws.Range("H49").Value = 6535.7334
$ws.Range("J49").Value = 6778.769
$ws.Range("L49").Value = 6778.769
$ws.Range("N49").Value = -7238.769
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H132").Value = 18299.445
$ws.Range("I132").Value = 15964
$ws.Range("J132").Value = 19467.166
$ws.Range("K132").Value = 47892
$ws.Range("L132").Value = 58401.49800000001
$ws.Range("M132").Value = -45362
$ws.Range("N132").Value = -63461.49800000001
